$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (and update the workbook's <sheet name="..."> entry) to the
# JETT template placeholder expression.
$ws.Name = "`${newSheetName}"

# Turn on "Normal" zoom display in the sheet view (zoomScaleNormal="100").
$excel.ActiveWindow.Zoom = 100

# Add the odd header / footer template placeholders.
$ws.PageSetup.LeftHeader   = "Header Left: `${numberList[0]}"
$ws.PageSetup.CenterHeader = "Header Center: `${numberList[1]}"
$ws.PageSetup.RightHeader  = "Header Right: `${numberList[2]}"
$ws.PageSetup.LeftFooter   = "Footer Left: `${numberList[3]}"
$ws.PageSetup.CenterFooter = "Footer Center: `${numberList[4]}"
$ws.PageSetup.RightFooter  = "Footer Right: `${numberList[5]}"
